$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 7 is currently empty (the sheet is sparse between row 6 and the TOTAL
# row 16). Fill it in directly instead of inserting, so the TOTAL row stays
# at row 16. 41064 is the serial date number for 2012-06-04 (matches the
# date serials already used in A2/A3).
$ws.Cells.Item(7, 1).Value = 41064
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = "ActionBar all version !"

# Match the formatting used by the rest of the table (date style, centered
# numeric style, left aligned text style) by copying the style from the row
# above.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(7, 1).Value = 41064
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = "ActionBar all version !"

$ws.Range("C8").Select()

# Touch the page setup so a <pageSetup/> element (portrait orientation) is
# emitted, matching the resave done by the original author's Excel.
$ws.PageSetup.Orientation = 1

$wb.Save()
